# Mark the following hubs as DONE ("[X] ") and fill in their feature
# columns (comment / first comment / community comment / original post)
# with "Y", matching the target diff:
#   - row 12: architecture      -> C,D,E,G = Y
#   - row 19: bnw                -> C,D,E,G = Y
#   - row 33: community_member   -> C,D,G   = Y   (E left blank)
#   - row 77: symmetry           -> C,D,G   = Y   (E left blank)
#   - row 87: wildlife           -> C,D,E,G = Y
# Row 30 (colorsplash) keeps its "[ ] " (not-done) status in column A,
# but still gets C,D,E,G filled in with "Y".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 12: architecture (mark done) ---
$ws.Range("A12").Value = "[X] "
$cols12 = @("C12", "D12", "E12", "G12")
foreach ($addr in $cols12) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = "Y"
}

# --- row 19: bnw (mark done) ---
$ws.Range("A19").Value = "[X] "
$cols19 = @("C19", "D19", "E19", "G19")
foreach ($addr in $cols19) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = "Y"
}

# --- row 30: colorsplash (stays not-done, only feature cells filled) ---
$cols30 = @("C30", "D30", "E30", "G30")
foreach ($addr in $cols30) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = "Y"
}

# --- row 33: community_member (mark done) ---
$ws.Range("A33").Value = "[X] "
$cols33 = @("C33", "D33", "G33")
foreach ($addr in $cols33) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = "Y"
}

# --- row 77: symmetry (mark done) ---
$ws.Range("A77").Value = "[X] "
$cols77 = @("C77", "D77", "G77")
foreach ($addr in $cols77) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = "Y"
}

# --- row 87: wildlife (mark done) ---
$ws.Range("A87").Value = "[X] "
$cols87 = @("C87", "D87", "E87", "G87")
foreach ($addr in $cols87) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = "Y"
}
